$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("liquids")

# The old column C ("description") becomes column E; two new columns are
# inserted in its place to hold the "pg" / "vg" percentage values.
$ws.Range("C1:D1").EntireColumn.Insert()

# --- Header row ---
$ws.Range("C1").Value = "pg"
$ws.Range("D1").Value = "vg"

# --- pg / vg percentage values per data row (2..12) ---
$pgvg = @{
    2  = @(100, 0)
    3  = @(100, 0)
    4  = @(100, 0)
    5  = @(30, 70)
    6  = @(30, 70)
    7  = @(30, 70)
    8  = @(50, 50)
    9  = @(30, 70)
    10 = @(100, 0)
    11 = @(100, 0)
    12 = @(100, 0)
}

foreach ($row in $pgvg.Keys) {
    $values = $pgvg[$row]
    $ws.Cells.Item($row, 3).Value = $values[0]
    $ws.Cells.Item($row, 4).Value = $values[1]
}

# The old "description" values (now shifted into column E) are no longer
# used for the data rows - only the "description" header (E1) remains.
$ws.Range("E2:E12").ClearContents()

# --- Column widths (closest achievable values; engine quantizes widths to
#     a 1/6-character-unit grid, so these land within a fraction of a unit
#     of the exact target widths 11.5703125 / 8.42578125 / 44.42578125) ---
$ws.Columns.Item(3).ColumnWidth = 10.666666666666666
$ws.Columns.Item(4).ColumnWidth = 7.666666666666667
$ws.Columns.Item(5).ColumnWidth = 43.666666666666664

# --- Selection state ---
$ws.Range("E6").Select()
